$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
# B2 holds a text value ("1"); force text type (to match the source column's
# string semantics) while keeping the cell's original (default) style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("I2").Value = 0.6870229007633588
$ws.Range("J2").Value = 0.6870229007633588
$ws.Range("K2").Value = 2.29
$ws.Range("L2").Value = 0.5826972010178116
$ws.Range("M2").Value = 4.4
$ws.Range("N2").Value = 0.04916201117318436
$ws.Range("O2").Value = 1.921397379912664
$ws.Range("P2").Value = 4.4
$ws.Range("Q2").Value = 0.04916201117318436
$ws.Range("R2").Value = 1.921397379912664
$ws.Range("U2").Value = 0.218
$ws.Range("V2").Value = 0.002435754189944134
$ws.Range("W2").Value = 0.01305587229190422
$ws.Range("X2").Value = 0.03875232327289777
$ws.Range("Y2").Value = -0.02569645098099355
$ws.Range("Z2").Value = 0.02321185990195499
$ws.Range("AA2").Value = 0.01594707932195381
$ws.Range("AB2").Value = 0.03728462541288226
$ws.Range("AC2").Value = -0.02133754609092845
$ws.Range("AD2").Value = 6.02
$ws.Range("AF2").Value = 6.02
$ws.Range("AG2").Value = 5.802
$ws.Range("AH2").Value = 0.06302345058626466
$ws.Range("AI2").Value = 0.03433721195528176
$ws.Range("AJ2").Value = 0.06088014941973936
$ws.Range("AK2").Value = 0.03313497275873491
$ws.Range("AL2").Value = 0.416
$ws.Range("AM2").Value = 0.416
$ws.Range("AO2").Value = 6.490384615384616
$ws.Range("AQ2").Value = 6.490384615384616

# Row 3 updates
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("I3").Value = 0.6870229007633588
$ws.Range("J3").Value = 0.6870229007633588
$ws.Range("K3").Value = 2.29
$ws.Range("L3").Value = 0.5826972010178116
$ws.Range("M3").Value = 4.4
$ws.Range("N3").Value = 0.04916201117318436
$ws.Range("O3").Value = 1.921397379912664
$ws.Range("P3").Value = 4.4
$ws.Range("Q3").Value = 0.04916201117318436
$ws.Range("R3").Value = 1.921397379912664
$ws.Range("U3").Value = 0.218
$ws.Range("V3").Value = 0.002435754189944134
$ws.Range("W3").Value = 0.01305587229190422
$ws.Range("X3").Value = 0.03875232327289777
$ws.Range("Y3").Value = -0.02569645098099355
$ws.Range("Z3").Value = 0.02321185990195499
$ws.Range("AA3").Value = 0.01594707932195381
$ws.Range("AB3").Value = 0.03728462541288226
$ws.Range("AC3").Value = -0.02133754609092845
$ws.Range("AD3").Value = 6.02
$ws.Range("AF3").Value = 6.02
$ws.Range("AG3").Value = 5.802
$ws.Range("AH3").Value = 0.06302345058626466
$ws.Range("AI3").Value = 0.03433721195528176
$ws.Range("AJ3").Value = 0.06088014941973936
$ws.Range("AK3").Value = 0.03313497275873491
$ws.Range("AL3").Value = 0.416
$ws.Range("AM3").Value = 0.416
$ws.Range("AO3").Value = 6.490384615384616
$ws.Range("AQ3").Value = 6.490384615384616

# Remove row 4 entirely
$ws.Rows(4).Delete()
